$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New header cells for row 1 (columns X, Y)
# ---------------------------------------------------------------------------
$ws.Range("X1").Value = "PriceChange"
$ws.Range("Y1").Value = "UpDown"

# ---------------------------------------------------------------------------
# Row 2 - update existing data row with new sentiment/trade values
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 42633.878854166665
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = "Buy"
$ws.Range("D2").Value = 34
$ws.Range("E2").Value = 6859
$ws.Range("F2").Value = 1266
$ws.Range("G2").Value = 63
$ws.Range("H2").Value = 35
$ws.Range("I2").Value = 94
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 11121
$ws.Range("L2").Value = 139
$ws.Range("M2").Value = 78
$ws.Range("N2").Value = 16
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = "Noun"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 1.76
$ws.Range("S2").Value = 0.111
$ws.Range("T2").Value = 5.45
$ws.Range("U2").Value = 4.84
$ws.Range("V2").Value = 2.2799999999999998
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = -1.6100000000000136
$ws.Range("Y2").Value = "Down"

# ---------------------------------------------------------------------------
# Row 3 - brand new data row
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = 42633.880277777775
$ws.Range("B3").Value = 11
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 48
$ws.Range("E3").Value = 6922
$ws.Range("F3").Value = 1309
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = 31
$ws.Range("I3").Value = 100
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 11942
$ws.Range("L3").Value = 146
$ws.Range("M3").Value = 66
$ws.Range("N3").Value = 20
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "Noun"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 1.76
$ws.Range("S3").Value = 0.111
$ws.Range("S3").NumberFormat = "0.00%"
$ws.Range("T3").Value = 5.45
$ws.Range("U3").Value = 4.84
$ws.Range("V3").Value = 2.2799999999999998
$ws.Range("W3").Value = 0

# ---------------------------------------------------------------------------
# Column C shrinks slightly now that "Buy"/"Down" replace the longer
# "Neutral" value (closest attainable width on this engine's quantized scale)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 6.67
